# Vwap.Calc.xlsx cleanup: rename the "index" column to "i" and switch it
# from a 1-based counter to a 0-based one; narrow column A to fit the
# shorter header.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("VWAP")

# Rename the header cell; the table (testdata7) column name follows the
# header cell automatically.
$ws.Range("A1").Value = "i"

# Re-index the data rows from 0 instead of 1 (rows 2..392 hold the 391
# data rows of the testdata7 table).
for ($r = 2; $r -le 392; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $cell.Value = $r - 2
}

# Column A only needs to fit a single-character header now.
$ws.Columns.Item(1).ColumnWidth = 3.2
